$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1644.6666
$ws.Range("I38").Value = 942.1111
$ws.Range("J38").Value = 2698.5
$ws.Range("K38").Value = 2826.3333
$ws.Range("L38").Value = 8095.5
$ws.Range("M38").Value = -2454.3333
$ws.Range("N38").Value = -8839.5
$ws.Range("H39").Value = 504.8095
$ws.Range("I39").Value = 529.6
$ws.Range("J39").Value = 9
$ws.Range("K39").Value = 1588.8
$ws.Range("L39").Value = 27
$ws.Range("M39").Value = -1292.8
$ws.Range("N39").Value = -619
$ws.Range("H99").Value = 2419.3845
$ws.Range("I99").Value = 1975.5
$ws.Range("J99").Value = 3899
$ws.Range("K99").Value = 5926.5
$ws.Range("L99").Value = 11697
$ws.Range("M99").Value = -4428.5
$ws.Range("N99").Value = -14693
$ws.Range("H101").Value = 3029
$ws.Range("I101").Value = 2760.1667
$ws.Range("J101").Value = 3175.6365
$ws.Range("K101").Value = 8280.500100000001
$ws.Range("L101").Value = 9526.9095
$ws.Range("M101").Value = -6658.500100000001
$ws.Range("N101").Value = -12770.9095
$ws.Range("H129").Value = 3671.75
$ws.Range("I129").Value = 4063.3333
$ws.Range("J129").Value = 2497
$ws.Range("K129").Value = 12189.9999
$ws.Range("L129").Value = 7491
$ws.Range("M129").Value = -7189.999899999999
$ws.Range("N129").Value = -17491
$ws.Range("H137").Value = 1857.5714
$ws.Range("I137").Value = 1936.6666
$ws.Range("J137").Value = 1383
$ws.Range("K137").Value = 5809.9998
$ws.Range("L137").Value = 4149
$ws.Range("M137").Value = -3259.9998
$ws.Range("N137").Value = -9249

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 12865.8
$ws.Range("I2").Value = 19554.666
$ws.Range("J2").Value = 2832.5
$ws.Range("K2").Value = 19554.666
$ws.Range("L2").Value = 2832.5
$ws.Range("M2").Value = -19441.666
$ws.Range("N2").Value = -3058.5
$ws.Range("H32").Value = 1360.0577
$ws.Range("I32").Value = 1140.4131
$ws.Range("J32").Value = 3044
$ws.Range("K32").Value = 1140.4131
$ws.Range("L32").Value = 3044
$ws.Range("M32").Value = -853.4131
$ws.Range("N32").Value = -3618
$ws.Range("H45").Value = 1351.8334
$ws.Range("I45").Value = 1037.3334
$ws.Range("K45").Value = 1037.3334
$ws.Range("M45").Value = -660.3334
$ws.Range("H116").Value = 12865.8
$ws.Range("I116").Value = 19554.666
$ws.Range("J116").Value = 2832.5
$ws.Range("K116").Value = 19554.666
$ws.Range("L116").Value = 2832.5
$ws.Range("M116").Value = -17260.666
$ws.Range("N116").Value = -7420.5
$ws.Range("H132").Value = 4015.0173
$ws.Range("I132").Value = 3986
$ws.Range("J132").Value = 4322.6
$ws.Range("K132").Value = 11958
$ws.Range("L132").Value = 12967.8
$ws.Range("M132").Value = -9428
$ws.Range("N132").Value = -18027.8

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 12865.8
$ws.Range("I3").Value = 19554.666
$ws.Range("J3").Value = 2832.5
$ws.Range("K3").Value = 19554.666
$ws.Range("L3").Value = 2832.5
$ws.Range("M3").Value = -19440.666
$ws.Range("N3").Value = -3060.5
$ws.Range("H86").Value = 143659.8
$ws.Range("I86").Value = 3766.6667
$ws.Range("J86").Value = 353499.5
$ws.Range("K86").Value = 3766.6667
$ws.Range("L86").Value = 353499.5
$ws.Range("M86").Value = -2643.6667
$ws.Range("N86").Value = -355745.5
$ws.Range("H89").Value = 143659.8
$ws.Range("I89").Value = 3766.6667
$ws.Range("J89").Value = 353499.5
$ws.Range("K89").Value = 18833.3335
$ws.Range("L89").Value = 1767497.5
$ws.Range("M89").Value = -13217.3335
$ws.Range("N89").Value = -1778729.5
$ws.Range("H105").Value = 2408.9412
$ws.Range("I105").Value = 1925.2142
$ws.Range("K105").Value = 1925.2142
$ws.Range("M105").Value = -178.2141999999999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 19665
$ws.Range("I56").Value = 14666.667
$ws.Range("J56").Value = 24663.334
$ws.Range("K56").Value = 14666.667
$ws.Range("L56").Value = 24663.334
$ws.Range("M56").Value = -13821.667
$ws.Range("N56").Value = -26353.334
$ws.Range("H102").Value = 26120.5
$ws.Range("J102").Value = 26120.5
$ws.Range("L102").Value = 26120.5
$ws.Range("N102").Value = -30988.5
$ws.Range("H122").Value = 1958.7
$ws.Range("I122").Value = 1916.4667
$ws.Range("J122").Value = 2085.4
$ws.Range("K122").Value = 5749.4001
$ws.Range("L122").Value = 6256.200000000001
$ws.Range("M122").Value = -3299.4001
$ws.Range("N122").Value = -11156.2
$ws.Range("H134").Value = 1702.7222
$ws.Range("I134").Value = 1702.7222
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5108.1666
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2573.1666
$ws.Range("N134").ClearContents()

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 2333.3333
$ws.Range("I9").Value = 2333.3333
$ws.Range("K9").Value = 6999.999899999999
$ws.Range("M9").Value = -6775.999899999999
$ws.Range("H68").Value = 1191.1321
$ws.Range("J68").Value = 1625.2858
$ws.Range("L68").Value = 4875.857400000001
$ws.Range("N68").Value = -6497.857400000001
$ws.Range("H71").Value = 1191.1321
$ws.Range("J71").Value = 1625.2858
$ws.Range("L71").Value = 14627.5722
$ws.Range("N71").Value = -22739.5722
$ws.Range("H107").Value = 1413.2273
$ws.Range("I107").Value = 182.66667
$ws.Range("J107").Value = 1607.5264
$ws.Range("K107").Value = 548.00001
$ws.Range("L107").Value = 4822.5792
$ws.Range("M107").Value = 1371.99999
$ws.Range("N107").Value = -8662.5792
$ws.Range("H109").Value = 2274.5715
$ws.Range("I109").Value = 1237
$ws.Range("J109").Value = 8500
$ws.Range("K109").Value = 3711
$ws.Range("L109").Value = 25500
$ws.Range("M109").Value = -2671
$ws.Range("N109").Value = -27580

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 11628499
$ws.Range("I14").Value = 18600800
$ws.Range("K14").Value = 18600800
$ws.Range("M14").Value = -18600632
$ws.Range("H113").Value = 2136.4167
$ws.Range("I113").Value = 2282.111
$ws.Range("J113").Value = 1699.3334
$ws.Range("K113").Value = 2282.111
$ws.Range("L113").Value = 1699.3334
$ws.Range("M113").Value = -112.1109999999999
$ws.Range("N113").Value = -6039.3334
$ws.Range("H132").Value = 2452.5278
$ws.Range("I132").Value = 2485.2
$ws.Range("J132").Value = 1309
$ws.Range("K132").Value = 7455.599999999999
$ws.Range("L132").Value = 3927
$ws.Range("M132").Value = -4925.599999999999
$ws.Range("N132").Value = -8987

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 1204499.8
$ws.Range("I12").Value = 2000166.6
$ws.Range("K12").Value = 2000166.6
$ws.Range("M12").Value = -1999996.6
$ws.Range("H99").Value = 20285
$ws.Range("J99").Value = 20285
$ws.Range("L99").Value = 20285
$ws.Range("N99").Value = -26275
$ws.Range("H121").Value = 40699
$ws.Range("J121").Value = 40699
$ws.Range("L121").Value = 40699
$ws.Range("N121").Value = -44193
$ws.Range("H132").Value = 2533.6667
$ws.Range("I132").Value = 2437.25
$ws.Range("J132").Value = 3305
$ws.Range("K132").Value = 7311.75
$ws.Range("L132").Value = 9915
$ws.Range("M132").Value = -4781.75
$ws.Range("N132").Value = -14975

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 24999
$ws.Range("J24").Value = 24999
$ws.Range("L24").Value = 24999
$ws.Range("N24").Value = -25459
$ws.Range("H109").Value = 49999
$ws.Range("J109").Value = 49999
$ws.Range("L109").Value = 49999
$ws.Range("N109").Value = -52773
$ws.Range("H132").Value = 2559.25
$ws.Range("I132").Value = 2559.25
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7677.75
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5147.75
$ws.Range("N132").ClearContents()

Write-Output "done"